$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Users")

# Insert 4 new rows before the current last row (row 27), pushing the
# existing last row (F00901 / 001) down to row 31.
$ws.Range("A27:C30").EntireRow.Insert()

$newIds = @("F02062", "F03153", "F00191", "F02729")
for ($i = 0; $i -lt $newIds.Length; $i++) {
    $r = 27 + $i
    $ws.Range("A$r").Value = $newIds[$i]
    $ws.Range("C$r").Value = "001"
    $ws.Range("C$r").NumberFormat = "@"
    $ws.Range("C$r").HorizontalAlignment = -4152
}

# Update the view so the newly-added rows are visible / selected, matching
# the saved workbook state.
$ws.Range("A31").Select()
$ws.Application.ActiveWindow.ScrollRow = 22
